# Fix heavy loadtest data corrected according to confluence Performance Test Data.xls

$wb = $excel.ActiveWorkbook

$heavy = $wb.Worksheets.Item("heavy")
$test  = $wb.Worksheets.Item("test")

# --- Data corrections on "heavy" sheet, row 6 (ConnectorLight) ---
$heavy.Range("E6").Value = 0
$heavy.Range("G6").Value = 0
# H6 used to be a formula (=100+K6); it becomes a hard-coded literal value.
$heavy.Range("H6").Value = 0
$heavy.Range("M6").Value = 0

# --- Data corrections on "heavy" sheet, row 7 (ConnectorMedium) ---
$heavy.Range("G7").Value = 0.25
$heavy.Range("M7").Value = 0.25

# --- View/selection changes ---
# "test" sheet keeps its own (now non-active) selection, updated to E6.
$test.Activate()
$test.Range("E6").Select()

# Active sheet moves from "test" back to "heavy", with a new selection/scroll position.
$heavy.Activate()
$heavy.Application.ActiveWindow.ScrollColumn = 4
$heavy.Range("D8").Select()
